{"js": "// Update the date line and the 25 multiplication problems in the table.\nconst replacements = [\n  [\"2025-10-05 Sunday\", \"2025-10-06 Monday\"],\n  [\"701\u00d75=\", \"872\u00d76=\"],\n  [\"970\u00d72=\", \"496\u00d76=\"],\n  [\"151\u00d79=\", \"811\u00d75=\"],\n  [\"991\u00d77=\", \"368\u00d76=\"],\n  [\"172\u00d74=\", \"159\u00d72=\"],\n  [\"943\u00d72=\", \"236\u00d74=\"],\n  [\"592\u00d76=\", \"421\u00d75=\"],\n  [\"398\u00d75=\", \"252\u00d76=\"],\n  [\"637\u00d77=\", \"471\u00d77=\"],\n  [\"665\u00d72=\", \"148\u00d76=\"],\n  [\"817\u00d77=\", \"234\u00d72=\"],\n  [\"950\u00d75=\", \"718\u00d72=\"],\n  [\"908\u00d74=\", \"848\u00d72=\"],\n  [\"731\u00d79=\", \"163\u00d73=\"],\n  [\"432\u00d76=\", \"894\u00d76=\"],\n  [\"999\u00d76=\", \"971\u00d79=\"],\n  [\"616\u00d77=\", \"591\u00d77=\"],\n  [\"936\u00d78=\", \"502\u00d79=\"],\n  [\"870\u00d78=\", \"780\u00d76=\"],\n  [\"232\u00d78=\", \"209\u00d78=\"],\n  [\"423\u00d74=\", \"579\u00d77=\"],\n  [\"766\u00d79=\", \"175\u00d76=\"],\n  [\"231\u00d76=\", \"499\u00d77=\"],\n  [\"937\u00d79=\", \"492\u00d72=\"],\n  [\"489\u00d77=\", \"800\u00d78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 multiplication problems in the table.\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-05 Sunday\", \"2025-10-06 Monday\"),\n    @(\"701\u00d75=\", \"872\u00d76=\"),\n    @(\"970\u00d72=\", \"496\u00d76=\"),\n    @(\"151\u00d79=\", \"811\u00d75=\"),\n    @(\"991\u00d77=\", \"368\u00d76=\"),\n    @(\"172\u00d74=\", \"159\u00d72=\"),\n    @(\"943\u00d72=\", \"236\u00d74=\"),\n    @(\"592\u00d76=\", \"421\u00d75=\"),\n    @(\"398\u00d75=\", \"252\u00d76=\"),\n    @(\"637\u00d77=\", \"471\u00d77=\"),\n    @(\"665\u00d72=\", \"148\u00d76=\"),\n    @(\"817\u00d77=\", \"234\u00d72=\"),\n    @(\"950\u00d75=\", \"718\u00d72=\"),\n    @(\"908\u00d74=\", \"848\u00d72=\"),\n    @(\"731\u00d79=\", \"163\u00d73=\"),\n    @(\"432\u00d76=\", \"894\u00d76=\"),\n    @(\"999\u00d76=\", \"971\u00d79=\"),\n    @(\"616\u00d77=\", \"591\u00d77=\"),\n    @(\"936\u00d78=\", \"502\u00d79=\"),\n    @(\"870\u00d78=\", \"780\u00d76=\"),\n    @(\"232\u00d78=\", \"209\u00d78=\"),\n    @(\"423\u00d74=\", \"579\u00d77=\"),\n    @(\"766\u00d79=\", \"175\u00d76=\"),\n    @(\"231\u00d76=\", \"499\u00d77=\"),\n    @(\"937\u00d79=\", \"492\u00d72=\"),\n    @(\"489\u00d77=\", \"800\u00d78=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, [ref]$newText, $wdReplaceOne) | Out-Null\n}\n"}
